$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(16,141,9,1,0.5,1,1,1,-0.5),
  @(16,142,-15,1,0.5,-1,0.7968627,1,-0.2968627),
  @(16,143,0,2,0.25,0,0,1,0.25),
  @(16,144,9,2,1,1,1,1,0),
  @(16,145,9,0.1,0,1,0.1950042,1,-0.1950042),
  @(16,146,-15,1,1,-1,0.9218271,1,0.07817289999999999),
  @(16,147,0,1,1,0,1,1,0),
  @(16,148,9,0.1,1,1,1,1,0),
  @(16,150,-15,1,0.5,-1,0.6484309,1,-0.1484309),
  @(16,151,9,1,0.25,1,0.2702366,1,-0.02023659999999999),
  @(16,152,-15,0.1,0,-1,0.1685711,1,-0.1685711),
  @(16,153,9,2,0.25,1,0.2539699,1,-0.003969900000000026),
  @(16,154,9,0.1,0.5,1,0.8578616999999999,1,-0.3578616999999999),
  @(16,155,-15,1,0.25,-1,0.2966696,1,-0.04666959999999998),
  @(16,156,9,1,1,1,0.8842947,1,0.1157053),
  @(16,157,-15,1,0,-1,0.2275373,1,-0.2275373),
  @(16,158,0,2,0.75,0,1,1,-0.25),
  @(16,160,9,0.1,0.5,1,0.5894653,1,-0.08946529999999997),
  @(16,161,-15,2,1,-1,1,1,0),
  @(16,162,9,2,0,1,0.2255039,1,-0.2255039),
  @(16,163,0,2,0.75,0,0.9452937,1,-0.1952937),
  @(16,164,9,0.1,0.75,1,0.9615609000000001,1,-0.2115609000000001),
  @(16,165,0,2,0.5,0,0.6972303,1,-0.1972303),
  @(16,166,9,2,0,1,0.319036,1,-0.319036),
  @(16,167,0,1,0.5,0,0.7378967,1,-0.2378967),
  @(16,168,0,1,0.75,0,1,1,-0.25),
  @(16,169,0,0.1,0,0,0.2051709,1,-0.2051709),
  @(16,170,9,2,0.5,1,0.9676601,1,-0.4676601),
  @(16,171,0,2,0.75,0,1,1,-0.25),
  @(16,172,0,2,1,0,0.9778268,1,0.0221732),
  @(16,173,-15,0.1,0.5,-1,0.8273624000000001,1,-0.3273624000000001),
  @(16,174,0,1,0.25,0,0.2011043,1,0.04889569999999999),
  @(16,175,0,0.1,1,0,0.9269941,1,0.07300589999999996),
  @(16,176,9,2,0.75,1,1,1,-0.25),
  @(16,177,-15,0.1,1,-1,1,1,0),
  @(16,178,-15,1,0.25,-1,0.2223705,1,0.0276295),
  @(16,179,-15,2,0.5,-1,0.2438036,1,0.2561964),
  @(16,180,0,0.1,0.5,0,0.5243997,1,-0.02439970000000002),
)

$startRow = 137
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
  }
}

Write-Host "Done"